# Weekly fruit/vegetable price update:
# A new observation row is inserted at row 198 (pushing the existing
# rows 198-209 down to 199-210), and the new row is populated with its
# own data (a new price report dated 44516).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, shifting rows 198:209 down to 199:210.
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly record.
$ws.Cells.Item(198, 1).Value  = 4
$ws.Cells.Item(198, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value  = "Los Lagos"
$ws.Cells.Item(198, 4).Value  = 44516
$ws.Cells.Item(198, 5).Value  = 10
$ws.Cells.Item(198, 6).Value  = 100114014
$ws.Cells.Item(198, 7).Value  = "Betarraga"
$ws.Cells.Item(198, 8).Value  = "Sin especificar"
$ws.Cells.Item(198, 9).Value  = "Primera"
$ws.Cells.Item(198, 10).Value = 1200
$ws.Cells.Item(198, 11).Value = 1000
$ws.Cells.Item(198, 12).Value = 1000
$ws.Cells.Item(198, 13).Value = 1000
$ws.Cells.Item(198, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 200
$ws.Cells.Item(198, 17).Value = 5
$ws.Cells.Item(198, 18).Value = "Hortaliza"
